# Add a new "N_cost" column (D) to the Tasks sheet and populate it with 1s,
# then leave the Tasks sheet as the selected/active sheet (it was previously
# the Supply sheet that was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# New header in D1
$ws.Range("D1").Value = "N_cost"

# Fill D2:D23 with the value 1
$ws.Range("D2:D23").Value = 1

# Make "Tasks" the active/selected sheet (it was "Supply" before the edit),
# and set the active cell to D2 as in the saved workbook.
$ws.Select()
$ws.Range("D2").Select()
